# ----------------------------------------------------------------------------
# Commit: "add some files, updates, store now all results in the same Databank"
#
# - Adds a new "TSFC Takeoff" column (E) with take-off TSFC values for the
#   first 28 engines, and refreshes their cruise-TSFC (column C) figures.
# - Folds 7 additional engines into the second databank block (which used to
#   occupy rows 30-52 with its own 0-based index) so everything now lives in
#   one contiguous table, rows 2-59, with the second block re-indexed 0..29
#   in rows 30-59. The take-off column is left blank for that block.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, formatted like the existing header cells (bold/border/center)
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("E1").Value = "TSFC Takeoff"

# Extend column A's bold/border/center style down through the 7 new rows
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A53:A59").PasteSpecial(-4122) | Out-Null

# Full data table, post-edit (row#, A, B, C, D, E) -- E = $null means "leave blank"
$data = @(
    ,@(2, 1, 1959, 27.73567457211094, 'B707-300', 22.47123130034523)
    ,@(3, 2, 1961, 26.63282981127525, 'B720-000', 21.14936708860759)
    ,@(4, 3, 1965, 22.09183986629342, 'DC9-10', 15.70655926352129)
    ,@(5, 4, 1966, 22.22154677600779, 'DC9-30', 15.8620253164557)
    ,@(6, 5, 1967, 23.0000762585199, 'B727-200/231A', 16.795166858458)
    ,@(7, 6, 1967, 22.54595806240666, 'B737-100/200', 16.25086306098964)
    ,@(8, 7, 1968, 23.0000762585199, 'DC9-40', 16.795166858458)
    ,@(9, 8, 1970, 17.09682773327364, 'DC10-10', 9.719562715765251)
    ,@(10, 9, 1970, 17.55084992131159, 'B747-200/300', 10.26375143843499)
    ,@(11, 10, 1970, 17.61575138020642, 'B747-100', 10.3415420023015)
    ,@(12, 11, 1972, 17.42114301159722, 'DC10-40', 10.10828538550058)
    ,@(13, 12, 1972, 17.74536228184551, 'DC10-30', 10.49689298043729)
    ,@(14, 13, 1973, 17.55084992131159, 'L1011-1/100/200', 10.26375143843499)
    ,@(15, 14, 1976, 21.70262312907501, 'DC9-50', 15.24004602991945)
    ,@(16, 15, 1979, 17.94025867468058, 'L1011-500', 10.73049482163406)
    ,@(17, 16, 1980, 20.59987437631462, 'MD80/DC9-80', 13.91829689298044)
    ,@(18, 17, 1983, 17.09682773327364, 'B767-200/ER', 9.719562715765251)
    ,@(19, 18, 1984, 16.90212335658914, 'A300-600', 9.486191024165707)
    ,@(20, 19, 1984, 17.42114301159722, 'B757-200', 10.10828538550058)
    ,@(21, 20, 1984, 20.72958128602899, 'B737-300', 14.07376294591485)
    ,@(22, 21, 1986, 16.90212335658914, 'A310-300', 9.486191024165707)
    ,@(23, 22, 1987, 16.96702481548398, 'B767-300/ER', 9.563981588032222)
    ,@(24, 23, 1988, 16.90212335658914, 'A320-100/200', 9.486191024165707)
    ,@(25, 24, 1988, 17.74536228184551, 'B737-400', 10.49689298043729)
    ,@(26, 25, 1989, 17.42114301159722, 'B747-400', 10.10828538550058)
    ,@(27, 26, 1990, 16.90212335658914, 'MD11', 9.486191024165707)
    ,@(28, 27, 1990, 17.94025867468058, 'B737-500/600', 10.73049482163406)
    ,@(29, 28, 1995, 15.86418005464829, 'B777', 8.242117376294594)
    ,@(30, 0, 1999, 17.42, '717-200                                   ', $null)
    ,@(31, 1, 1964, 22.52007039498443, '727-100                                   ', $null)
    ,@(32, 2, 1968, 22.44014078996885, '737-200C                                   ', $null)
    ,@(33, 3, 1997, 17.28773126870821, '737-700/700LR/Max 7                       ', $null)
    ,@(34, 4, 1998, 17.63566548921198, '737-800                                 ', $null)
    ,@(35, 5, 2001, 17.74588377645939, '737-900                                  ', $null)
    ,@(36, 6, 2007, 17.63444277935031, '737-900ER                                 ', $null)
    ,@(37, 7, 1999, 17.278741900219, '757-300                                   ', $null)
    ,@(38, 8, 2000, 17.00474996911675, '767-400/ER                                ', $null)
    ,@(39, 9, 1997, 16.29765369197297, '777-300/300ER/333ER                       ', $null)
    ,@(40, 10, 2018, 15.33698632247987, '787-10 Dreamliner                          ', $null)
    ,@(41, 11, 1996, 16.7661071789652, 'A319                        ', $null)
    ,@(42, 12, 2016, 14.9163310818809, 'A320-200n                        ', $null)
    ,@(43, 13, 2007, 14.99919751506516, 'A321-200n                        ', $null)
    ,@(44, 14, 1994, 17.42351819037375, 'A321/Lr                        ', $null)
    ,@(45, 15, 1998, 16.9928785194305, 'A330-200                    ', $null)
    ,@(46, 16, 1994, 16.5453498438535, 'A330-300/333                    ', $null)
    ,@(47, 17, 2018, 15.30993156577482, 'A330-900                         ', $null)
    ,@(48, 18, 2017, 15.73636826305148, 'B737 Max 800                               ', $null)
    ,@(49, 19, 2018, 15.73636826305148, 'B737 Max 900   ', $null)
    ,@(50, 20, 2011, 15.33698632247987, 'B787-800 Dreamliner                              ', $null)
    ,@(51, 21, 2014, 15.33698632247987, 'B787-900 Dreamliner                              ', $null)
    ,@(52, 22, 2001, 18.13689908825537, 'CRJ 900                                 ', $null)
    ,@(53, 23, 2004, 17.58758235398228, 'Embraer 190                                      ', $null)
    ,@(54, 24, 2005, 18.19612006630815, 'Embraer ERJ-175                                  ', $null)
    ,@(55, 25, 2001, 18.12812519816411, 'Embraer-140                                      ', $null)
    ,@(56, 26, 1996, 18.13073907793062, 'Embraer-145', $null)
    ,@(57, 27, 1995, 16.88813484187854, 'MD-90                         ', $null)
    ,@(58, 28, 1992, 19.49, 'RJ-200ER /RJ-440                   ', $null)
    ,@(59, 29, 1999, 17.95993678687839, 'RJ-700                                  ', $null)
)

foreach ($row in $data) {
    $r  = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    if ($row[5] -ne $null) {
        $ws.Cells.Item($r, 5).Value = $row[5]
    }
}
